# Completed Documentation of Level 1
# - Update the title from "Website D" to "Website B"
# - Mark the edit location with Word's "_GoBack" bookmark (the bookmark
#   Word itself drops at the last edit position), placed right after the
#   title run, inside the title paragraph.

$d = $word.ActiveDocument

# 1. Update the title text, editing the Range.Text in place so the
#    existing run (and its rsid attributes) are preserved rather than
#    replaced wholesale.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Text = "Documentation : Website B"

# 2. Insert the "_GoBack" bookmark immediately after the title text,
#    still inside the title paragraph (collapsed / zero-length bookmark).
#    A temporary marker character is used so the insertion point sits
#    unambiguously at the end of the run's text (rather than gravitating
#    to the following paragraph), then the marker is removed again.
$titlePara = $d.Paragraphs(1)
$titleEnd = $titlePara.Range.End - 1

$marker = $d.Range($titleEnd, $titleEnd)
$marker.InsertAfter("~")

$bookmarkRange = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$markerRange = $d.Range($titleEnd, $titleEnd + 1)
$markerRange.Delete()

Write-Output "Title updated and _GoBack bookmark inserted."
